$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 8 (pushes existing "extr1".."extr8" rows down
# by two, from rows 8-15 to rows 10-17) and restore the border formatting
# that the engine drops on the freshly inserted rows so the style matches
# the other data rows (cellXfs index used by column A, s="1").
$ws.Rows("8:9").Insert()
$ws.Range("A8:A9").Borders.LineStyle = 1

# The "A" column is just a running index (row number - 2); the insert
# shifted the old values down with the rows, so renumber column A back
# into a single unbroken sequence 0..15 for rows 2..17.
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10
$ws.Range("A13").Value = 11
$ws.Range("A14").Value = 12
$ws.Range("A15").Value = 13
$ws.Range("A16").Value = 14
$ws.Range("A17").Value = 15

# Populate the two new rows with the "line7" / "line8" data (B10:E17 keep
# the "extr1".."extr8" data that the insert already shifted down intact).
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true
